$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "'46.097.15"
$ws.Cells.Item(2, 4).Style = "Normal"
$ws.Cells.Item(2, 5).Value = '  +1.43%  '
$ws.Cells.Item(3, 4).Value = "'2.588.75"
$ws.Cells.Item(3, 4).Style = "Normal"
$ws.Cells.Item(3, 5).Value = '  +7.19%  '
$ws.Cells.Item(4, 5).Value = '  +0.10%  '
$ws.Cells.Item(5, 4).Value = "'306.58"
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value = '  +4.10%  '
$ws.Cells.Item(6, 4).Value = "'99.16"
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Cells.Item(6, 5).Value = '  +4.71%  '
$ws.Cells.Item(7, 4).Value = "'0.601"
$ws.Cells.Item(7, 4).Style = "Normal"
$ws.Cells.Item(7, 5).Value = '  +6.50%  '
$ws.Cells.Item(8, 5).Value = '  +0.14%  '
$ws.Cells.Item(9, 4).Value = "'0.578"
$ws.Cells.Item(9, 4).Style = "Normal"
$ws.Cells.Item(9, 5).Value = '  +15.29%  '
$ws.Cells.Item(10, 4).Value = "'39.22"
$ws.Cells.Item(10, 4).Style = "Normal"
$ws.Cells.Item(10, 5).Value = '  +12.16%  '
$ws.Cells.Item(11, 4).Value = "'54.31"
$ws.Cells.Item(11, 4).Style = "Normal"
$ws.Cells.Item(11, 5).Value = '  +1.47%  '
$ws.Cells.Item(12, 4).Value = "'0.0841"
$ws.Cells.Item(12, 4).Style = "Normal"
$ws.Cells.Item(12, 5).Value = '  +7.69%  '
$ws.Cells.Item(13, 4).Value = "'8.18"
$ws.Cells.Item(13, 4).Style = "Normal"
$ws.Cells.Item(13, 5).Value = '  +16.02%  '
$ws.Cells.Item(14, 4).Value = "'2.990.19"
$ws.Cells.Item(14, 4).Style = "Normal"
$ws.Cells.Item(14, 5).Value = '  +7.15%  '
$ws.Cells.Item(15, 5).Value = '  +1.46%  '
$ws.Cells.Item(16, 4).Value = "'2.611.11"
$ws.Cells.Item(16, 4).Style = "Normal"
$ws.Cells.Item(16, 5).Value = '  +7.77%  '
$ws.Cells.Item(17, 4).Value = "'0.914"
$ws.Cells.Item(17, 4).Style = "Normal"
$ws.Cells.Item(17, 5).Value = '  +9.29%  '
$ws.Cells.Item(18, 4).Value = "'14.91"
$ws.Cells.Item(18, 4).Style = "Normal"
$ws.Cells.Item(18, 5).Value = '  +5.93%  '
$ws.Cells.Item(19, 4).Value = "'46.362.80"
$ws.Cells.Item(19, 4).Style = "Normal"
$ws.Cells.Item(19, 5).Value = '  +2.20%  '
$ws.Cells.Item(20, 5).Value = '  +7.00%  '
$ws.Cells.Item(21, 5).Value = '  +4.75%  '
$ws.Cells.Item(22, 4).Value = "'6.66"
$ws.Cells.Item(22, 4).Style = "Normal"
$ws.Cells.Item(22, 5).Value = '  +8.14%  '
$ws.Cells.Item(23, 4).Value = "'71.70"
$ws.Cells.Item(23, 4).Style = "Normal"
$ws.Cells.Item(23, 5).Value = '  +6.90%  '
$ws.Cells.Item(24, 4).Value = "'270.86"
$ws.Cells.Item(24, 4).Style = "Normal"
$ws.Cells.Item(24, 5).Value = '  +12.40%  '
$ws.Cells.Item(25, 5).Value = '  +8.27%  '
$ws.Cells.Item(26, 4).Value = "'30.18"
$ws.Cells.Item(26, 4).Style = "Normal"
$ws.Cells.Item(26, 5).Value = '  +42.77%  '
$ws.Cells.Item(27, 4).Value = "'2.16"
$ws.Cells.Item(27, 4).Style = "Normal"
$ws.Cells.Item(27, 5).Value = '  +11.36%  '
$ws.Cells.Item(28, 5).Value = '  +0.04%  '
$ws.Cells.Item(29, 4).Value = "'4.03"
$ws.Cells.Item(29, 4).Style = "Normal"
$ws.Cells.Item(29, 5).Value = '  +0.28%  '
$ws.Cells.Item(30, 4).Value = "'10.52"
$ws.Cells.Item(30, 4).Style = "Normal"
$ws.Cells.Item(30, 5).Value = '  +8.55%  '
$ws.Cells.Item(31, 5).Value = '  +3.75%  '
$ws.Cells.Item(32, 4).Value = "'39.26"
$ws.Cells.Item(32, 4).Style = "Normal"
$ws.Cells.Item(32, 5).Value = '  +1.13%  '
$ws.Cells.Item(33, 5).Value = '  +13.21%  '
$ws.Cells.Item(34, 4).Value = "'3.62"
$ws.Cells.Item(34, 4).Style = "Normal"
$ws.Cells.Item(34, 5).Value = '  -3.61%  '
$ws.Cells.Item(35, 5).Value = '  +3.47%  '
$ws.Cells.Item(36, 4).Value = "'0.0836"
$ws.Cells.Item(36, 4).Style = "Normal"
$ws.Cells.Item(36, 5).Value = '  +9.37%  '
$ws.Cells.Item(37, 5).Value = '  +10.70%  '
$ws.Cells.Item(38, 4).Value = "'149.78"
$ws.Cells.Item(38, 4).Style = "Normal"
$ws.Cells.Item(38, 5).Value = '  +0.67%  '
$ws.Cells.Item(39, 5).Value = '  +8.48%  '
$ws.Cells.Item(40, 5).Value = '  +5.37%  '
$ws.Cells.Item(41, 4).Value = "'23.14"
$ws.Cells.Item(41, 4).Style = "Normal"
$ws.Cells.Item(41, 5).Value = '  +43.95%  '
$ws.Cells.Item(42, 4).Value = "'16.09"
$ws.Cells.Item(42, 4).Style = "Normal"
$ws.Cells.Item(42, 5).Value = '  +8.45%  '
$ws.Cells.Item(43, 4).Value = "'0.0327"
$ws.Cells.Item(43, 4).Style = "Normal"
$ws.Cells.Item(43, 5).Value = '  +10.00%  '
$ws.Cells.Item(44, 4).Value = "'3.59"
$ws.Cells.Item(44, 4).Style = "Normal"
$ws.Cells.Item(44, 5).Value = '  +11.60%  '
$ws.Cells.Item(45, 4).Value = "'4.07"
$ws.Cells.Item(45, 4).Style = "Normal"
$ws.Cells.Item(45, 5).Value = '  +7.80%  '
$ws.Cells.Item(46, 4).Value = "'2.140.89"
$ws.Cells.Item(46, 4).Style = "Normal"
$ws.Cells.Item(46, 5).Value = '  +6.59%  '
$ws.Cells.Item(47, 4).Value = "'0.999"
$ws.Cells.Item(47, 4).Style = "Normal"
$ws.Cells.Item(47, 5).Value = '  -0.05%  '
$ws.Cells.Item(48, 4).Value = "'93.29"
$ws.Cells.Item(48, 4).Style = "Normal"
$ws.Cells.Item(48, 5).Value = '  +4.66%  '
$ws.Cells.Item(49, 4).Value = "'9.70"
$ws.Cells.Item(49, 4).Style = "Normal"
$ws.Cells.Item(49, 5).Value = '  +13.95%  '
$ws.Cells.Item(50, 4).Value = "'1.77"
$ws.Cells.Item(50, 4).Style = "Normal"
$ws.Cells.Item(50, 5).Value = '  -0.20%  '
$ws.Cells.Item(51, 4).Value = "'108.52"
$ws.Cells.Item(51, 4).Style = "Normal"
$ws.Cells.Item(51, 5).Value = '  +7.62%  '
